$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two columns before column A
$ws.Range("A:B").Insert()

# Insert one column after the (now shifted) last column (old Z -> now AB), so insert before AC
$ws.Range("AC:AC").Insert()

# Copy formatting from the neighboring header cell (C1, same style as all the
# other general headers) onto the two new leading columns.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats

# Copy formatting from Z1 (date-style header) onto the new trailing column.
$ws.Range("Z1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set header values for the two new leading columns
# (set B1 first so shared-string index order matches: 26=Número de Carga, 27=Validación por Supervisión)
$ws.Range("B1").Value = "Número de Carga"
$ws.Range("A1").Value = "Validación por Supervisión"

# Set header value for the new trailing column
$ws.Range("AC1").Value = "Observaciones"

# Set column widths (closest reachable values; ColumnWidth is rounded to
# whole-pixel granularity internally, same as real Excel does)
$ws.Range("A:B").ColumnWidth = 12.67
$ws.Range("AC:AC").ColumnWidth = 13

# Reset selection back to the default top-left cell (matches the workbook
# being resaved with no special selection left on it)
$ws.Range("A1").Select() | Out-Null
